$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Constants ---
$xlPasteFormats = -4122

# -----------------------------------------------------------------
# renameRelation block (row 45-46): split the old single "relation"
# parameter into two dedicated parameters.
# -----------------------------------------------------------------
$ws.Cells.Item(45, 5).Value2 = "relationToRename"   # E45: relation -> relationToRename
$ws.Cells.Item(46, 5).Value2 = "newRelationName"    # E46: relationname -> newRelationName

# -----------------------------------------------------------------
# ChangeRelationCardinality block (row 47-48): add the missing
# function-name cell, rename the parameter, and extend row 48/47
# with the same row "skeleton" (borders) used by the sibling blocks.
# -----------------------------------------------------------------

# New trailing border-only cell H47 - copy formatting from H45 (style 18)
$ws.Range("H45").Copy() | Out-Null
$ws.Range("H47").PasteSpecial($xlPasteFormats) | Out-Null

# Row 48 gains the full set of bottom-border cells (A:D, G) like row 46,
# and E48/F48 switch from the "row45/47" border (style 2) to the
# "row46" border (style 1). Re-use row 46 as the formatting template.
$ws.Range("A46").Copy() | Out-Null
$ws.Range("A48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B46").Copy() | Out-Null
$ws.Range("B48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("C46").Copy() | Out-Null
$ws.Range("C48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D46").Copy() | Out-Null
$ws.Range("D48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("E46").Copy() | Out-Null
$ws.Range("E48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("F46").Copy() | Out-Null
$ws.Range("F48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("G46").Copy() | Out-Null
$ws.Range("G48").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("H46").Copy() | Out-Null
$ws.Range("H48").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# Parameter name/value updates - order matters for shared-string append
# order: newCardinality (row 48) is introduced before the new function
# name changeCardinalityInRelation (row 47, column A).
$ws.Cells.Item(48, 5).Value2 = "newCardinality"   # E48: cardinality (newCardinality?) -> newCardinality
$ws.Cells.Item(48, 6).Value2 = "Cardinality"      # F48: unchanged text, re-affirm after paste
$ws.Cells.Item(47, 1).Value2 = "changeCardinalityInRelation"   # A47: new function name cell

# -----------------------------------------------------------------
# Scroll position / selection bookkeeping (view-state only change)
# -----------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 1
$ws.Range("C40").Select() | Out-Null
